# Auto-generated edit script applying the Golem_Profits crafting-profit refresh
# described by the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 78  # H38: 398 -> 78
$ws.Cells.Item(38, 9).Value = 78  # I38: 398 -> 78
$ws.Cells.Item(38, 11).Value = 234  # K38: 1194 -> 234
$ws.Cells.Item(38, 13).Value = 138  # M38: -822 -> 138
$ws.Cells.Item(43, 8).Value = 4570.8  # H43: 4963.5 -> 4570.8
$ws.Cells.Item(43, 9).Value = 3450  # I43: 3900 -> 3450
$ws.Cells.Item(43, 11).Value = 3450  # K43: 3900 -> 3450
$ws.Cells.Item(43, 13).Value = -3381  # M43: -3831 -> -3381
$ws.Cells.Item(58, 8).Value = 533.625  # H58: 1012.25 -> 533.625
$ws.Cells.Item(58, 9).Value = 538.4286  # I58: 1012.25 -> 538.4286
$ws.Cells.Item(58, 10).Value = 500  # J58: 0 -> 500
$ws.Cells.Item(58, 11).Value = 1615.2858  # K58: 3036.75 -> 1615.2858
$ws.Cells.Item(58, 12).Value = 1500  # L58: 0 -> 1500
$ws.Cells.Item(58, 13).Value = -1465.2858  # M58: -2886.75 -> -1465.2858
$ws.Cells.Item(58, 14).Value = -1800  # N58: None -> -1800
$ws.Cells.Item(86, 8).Value = 6889.222  # H86: 10625.5 -> 6889.222
$ws.Cells.Item(86, 9).Value = 3000  # I86: 3500 -> 3000
$ws.Cells.Item(86, 10).Value = 8833.833000000001  # J86: 13000.667 -> 8833.833000000001
$ws.Cells.Item(86, 11).Value = 3000  # K86: 3500 -> 3000
$ws.Cells.Item(86, 12).Value = 8833.833000000001  # L86: 13000.667 -> 8833.833000000001
$ws.Cells.Item(86, 13).Value = -1877  # M86: -2377 -> -1877
$ws.Cells.Item(86, 14).Value = -11079.833  # N86: -15246.667 -> -11079.833
$ws.Cells.Item(87, 8).Value = 110000  # H87: 99000 -> 110000
$ws.Cells.Item(87, 10).Value = 110000  # J87: 99000 -> 110000
$ws.Cells.Item(87, 12).Value = 110000  # L87: 99000 -> 110000
$ws.Cells.Item(87, 14).Value = -112496  # N87: -101496 -> -112496
$ws.Cells.Item(89, 8).Value = 6889.222  # H89: 10625.5 -> 6889.222
$ws.Cells.Item(89, 9).Value = 3000  # I89: 3500 -> 3000
$ws.Cells.Item(89, 10).Value = 8833.833000000001  # J89: 13000.667 -> 8833.833000000001
$ws.Cells.Item(89, 11).Value = 15000  # K89: 17500 -> 15000
$ws.Cells.Item(89, 12).Value = 44169.165  # L89: 65003.335 -> 44169.165
$ws.Cells.Item(89, 13).Value = -9384  # M89: -11884 -> -9384
$ws.Cells.Item(89, 14).Value = -55401.165  # N89: -76235.33499999999 -> -55401.165
$ws.Cells.Item(90, 8).Value = 110000  # H90: 99000 -> 110000
$ws.Cells.Item(90, 10).Value = 110000  # J90: 99000 -> 110000
$ws.Cells.Item(90, 12).Value = 330000  # L90: 297000 -> 330000
$ws.Cells.Item(90, 14).Value = -342480  # N90: -309480 -> -342480
$ws.Cells.Item(135, 8).Value = 1388  # H135: 1495 -> 1388
$ws.Cells.Item(135, 9).Value = 1333.3334  # I135: 1500 -> 1333.3334
$ws.Cells.Item(135, 10).Value = 1470  # J135: 1485 -> 1470
$ws.Cells.Item(135, 11).Value = 12000.0006  # K135: 13500 -> 12000.0006
$ws.Cells.Item(135, 12).Value = 13230  # L135: 13365 -> 13230
$ws.Cells.Item(135, 13).Value = -9465.000599999999  # M135: -10965 -> -9465.000599999999
$ws.Cells.Item(135, 14).Value = -18300  # N135: -18435 -> -18300
$ws.Cells.Item(137, 8).Value = 2621  # H137: 2663 -> 2621
$ws.Cells.Item(137, 9).Value = 2494.6667  # I137: 2494.5 -> 2494.6667
$ws.Cells.Item(137, 11).Value = 7484.000100000001  # K137: 7483.5 -> 7484.000100000001
$ws.Cells.Item(137, 13).Value = -4934.000100000001  # M137: -4933.5 -> -4934.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 2111  # H35: 3016 -> 2111
$ws.Cells.Item(35, 9).Value = 2111  # I35: 3016 -> 2111
$ws.Cells.Item(35, 11).Value = 2111  # K35: 3016 -> 2111
$ws.Cells.Item(35, 13).Value = -1705  # M35: -2610 -> -1705
$ws.Cells.Item(61, 8).Value = 1270.3334  # H61: 1251.5 -> 1270.3334
$ws.Cells.Item(61, 9).Value = 1306  # I61: 1269 -> 1306
$ws.Cells.Item(61, 11).Value = 1306  # K61: 1269 -> 1306
$ws.Cells.Item(61, 13).Value = -1094  # M61: -1057 -> -1094
$ws.Cells.Item(88, 8).Value = 2941.2856  # H88: 3012 -> 2941.2856
$ws.Cells.Item(88, 10).Value = 3838  # J88: 3937 -> 3838
$ws.Cells.Item(88, 12).Value = 3838  # L88: 3937 -> 3838
$ws.Cells.Item(88, 14).Value = -4650  # N88: -4749 -> -4650
$ws.Cells.Item(91, 8).Value = 2941.2856  # H91: 3012 -> 2941.2856
$ws.Cells.Item(91, 10).Value = 3838  # J91: 3937 -> 3838
$ws.Cells.Item(91, 12).Value = 3838  # L91: 3937 -> 3838
$ws.Cells.Item(91, 14).Value = -6646  # N91: -6745 -> -6646
$ws.Cells.Item(97, 8).Value = 1934.7  # H97: 2069.7 -> 1934.7
$ws.Cells.Item(97, 9).Value = 1499.625  # I97: 1585.2858 -> 1499.625
$ws.Cells.Item(97, 10).Value = 3675  # J97: 3200 -> 3675
$ws.Cells.Item(97, 11).Value = 1499.625  # K97: 1585.2858 -> 1499.625
$ws.Cells.Item(97, 12).Value = 3675  # L97: 3200 -> 3675
$ws.Cells.Item(97, 13).Value = -1003.625  # M97: -1089.2858 -> -1003.625
$ws.Cells.Item(97, 14).Value = -4667  # N97: -4192 -> -4667
$ws.Cells.Item(136, 8).Value = 1270.3334  # H136: 1251.5 -> 1270.3334
$ws.Cells.Item(136, 9).Value = 1306  # I136: 1269 -> 1306
$ws.Cells.Item(136, 11).Value = 3918  # K136: 3807 -> 3918
$ws.Cells.Item(136, 13).Value = -1368  # M136: -1257 -> -1368

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2454.5454  # H94: 2363.5454 -> 2454.5454
$ws.Cells.Item(94, 9).Value = 2000  # I94: 1833.1666 -> 2000
$ws.Cells.Item(94, 11).Value = 2000  # K94: 1833.1666 -> 2000
$ws.Cells.Item(94, 13).Value = -1549  # M94: -1382.1666 -> -1549

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 99000  # H68: 0 -> 99000
$ws.Cells.Item(68, 10).Value = 99000  # J68: 0 -> 99000
$ws.Cells.Item(68, 12).Value = 99000  # L68: 0 -> 99000
$ws.Cells.Item(68, 14).Value = -100498  # N68: None -> -100498
$ws.Cells.Item(71, 8).Value = 99000  # H71: 0 -> 99000
$ws.Cells.Item(71, 10).Value = 99000  # J71: 0 -> 99000
$ws.Cells.Item(71, 12).Value = 297000  # L71: 0 -> 297000
$ws.Cells.Item(71, 14).Value = -304488  # N71: None -> -304488
$ws.Cells.Item(125, 8).Value = 18845.666  # H125: 19089.666 -> 18845.666
$ws.Cells.Item(125, 10).Value = 18845.666  # J125: 19089.666 -> 18845.666
$ws.Cells.Item(125, 12).Value = 18845.666  # L125: 19089.666 -> 18845.666
$ws.Cells.Item(125, 14).Value = -23765.666  # N125: -24009.666 -> -23765.666
$ws.Cells.Item(132, 8).Value = 0  # H132: 1500 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 1500 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 4500 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: -1970 -> (removed)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 1750  # H59: 0 -> 1750
$ws.Cells.Item(59, 10).Value = 1750  # J59: 0 -> 1750
$ws.Cells.Item(59, 12).Value = 5250  # L59: 0 -> 5250
$ws.Cells.Item(59, 14).Value = -6330  # N59: None -> -6330
$ws.Cells.Item(99, 8).Value = 3012.5  # H99: 4217 -> 3012.5
$ws.Cells.Item(99, 9).Value = 3012.5  # I99: 4217 -> 3012.5
$ws.Cells.Item(99, 11).Value = 9037.5  # K99: 12651 -> 9037.5
$ws.Cells.Item(99, 13).Value = -6791.5  # M99: -10405 -> -6791.5
$ws.Cells.Item(109, 8).Value = 1007  # H109: 907 -> 1007
$ws.Cells.Item(109, 9).Value = 1114  # I109: 909.3333 -> 1114
$ws.Cells.Item(109, 11).Value = 3342  # K109: 2727.9999 -> 3342
$ws.Cells.Item(109, 13).Value = -2302  # M109: -1687.9999 -> -2302
$ws.Cells.Item(111, 8).Value = 349  # H111: 0 -> 349
$ws.Cells.Item(111, 9).Value = 349  # I111: 0 -> 349
$ws.Cells.Item(111, 11).Value = 1047  # K111: 0 -> 1047
$ws.Cells.Item(111, 13).Value = 2020  # M111: None -> 2020
$ws.Cells.Item(114, 8).Value = 2165  # H114: 2499.75 -> 2165
$ws.Cells.Item(114, 9).Value = 1750  # I114: 1666.3334 -> 1750
$ws.Cells.Item(114, 10).Value = 2995  # J114: 5000 -> 2995
$ws.Cells.Item(114, 11).Value = 5250  # K114: 4999.0002 -> 5250
$ws.Cells.Item(114, 12).Value = 8985  # L114: 15000 -> 8985
$ws.Cells.Item(114, 13).Value = -1996  # M114: -1745.0002 -> -1996
$ws.Cells.Item(114, 14).Value = -15493  # N114: -21508 -> -15493
$ws.Cells.Item(129, 8).Value = 2689.5  # H129: 2347.8572 -> 2689.5
$ws.Cells.Item(129, 9).Value = 368  # I129: 350.5 -> 368
$ws.Cells.Item(129, 11).Value = 1104  # K129: 1051.5 -> 1104
$ws.Cells.Item(129, 13).Value = 3896  # M129: 3948.5 -> 3896

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 2174.875  # H70: 3041.6 -> 2174.875
$ws.Cells.Item(70, 9).Value = 2174.875  # I70: 2833 -> 2174.875
$ws.Cells.Item(70, 10).Value = 0  # J70: 3354.5 -> 0
$ws.Cells.Item(70, 11).Value = 2174.875  # K70: 2833 -> 2174.875
$ws.Cells.Item(70, 12).Value = 0  # L70: 3354.5 -> 0
$ws.Cells.Item(70, 13).Value = -1904.875  # M70: -2563 -> -1904.875
$ws.Cells.Item(70, 14).ClearContents()  # N70: -3894.5 -> (removed)
$ws.Cells.Item(73, 8).Value = 2174.875  # H73: 3041.6 -> 2174.875
$ws.Cells.Item(73, 9).Value = 2174.875  # I73: 2833 -> 2174.875
$ws.Cells.Item(73, 10).Value = 0  # J73: 3354.5 -> 0
$ws.Cells.Item(73, 11).Value = 2174.875  # K73: 2833 -> 2174.875
$ws.Cells.Item(73, 12).Value = 0  # L73: 3354.5 -> 0
$ws.Cells.Item(73, 13).Value = -1238.875  # M73: -1897 -> -1238.875
$ws.Cells.Item(73, 14).ClearContents()  # N73: -5226.5 -> (removed)
$ws.Cells.Item(80, 8).Value = 780  # H80: 0 -> 780
$ws.Cells.Item(80, 9).Value = 945  # I80: 0 -> 945
$ws.Cells.Item(80, 10).Value = 450  # J80: 0 -> 450
$ws.Cells.Item(80, 11).Value = 945  # K80: 0 -> 945
$ws.Cells.Item(80, 12).Value = 450  # L80: 0 -> 450
$ws.Cells.Item(80, 13).Value = 53  # M80: None -> 53
$ws.Cells.Item(80, 14).Value = -2446  # N80: None -> -2446
$ws.Cells.Item(83, 8).Value = 780  # H83: 0 -> 780
$ws.Cells.Item(83, 9).Value = 945  # I83: 0 -> 945
$ws.Cells.Item(83, 10).Value = 450  # J83: 0 -> 450
$ws.Cells.Item(83, 11).Value = 4725  # K83: 0 -> 4725
$ws.Cells.Item(83, 12).Value = 2250  # L83: 0 -> 2250
$ws.Cells.Item(83, 13).Value = 267  # M83: None -> 267
$ws.Cells.Item(83, 14).Value = -12234  # N83: None -> -12234
$ws.Cells.Item(102, 8).Value = 45867.777  # H102: 45901 -> 45867.777
$ws.Cells.Item(102, 9).Value = 51401.25  # I102: 51438.625 -> 51401.25
$ws.Cells.Item(102, 11).Value = 51401.25  # K102: 51438.625 -> 51401.25
$ws.Cells.Item(102, 13).Value = -49779.25  # M102: -49816.625 -> -49779.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 602.6667  # H9: 750 -> 602.6667
$ws.Cells.Item(9, 9).Value = 777.5  # I9: 500 -> 777.5
$ws.Cells.Item(9, 10).Value = 253  # J9: 1000 -> 253
$ws.Cells.Item(9, 11).Value = 777.5  # K9: 500 -> 777.5
$ws.Cells.Item(9, 12).Value = 253  # L9: 1000 -> 253
$ws.Cells.Item(9, 13).Value = -553.5  # M9: -276 -> -553.5
$ws.Cells.Item(9, 14).Value = -701  # N9: -1448 -> -701
$ws.Cells.Item(20, 8).Value = 1249.5  # H20: 1333.3334 -> 1249.5
$ws.Cells.Item(20, 10).Value = 1249.5  # J20: 1333.3334 -> 1249.5
$ws.Cells.Item(20, 12).Value = 1249.5  # L20: 1333.3334 -> 1249.5
$ws.Cells.Item(20, 14).Value = -1701.5  # N20: -1785.3334 -> -1701.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1325  # H122: 1400 -> 1325
$ws.Cells.Item(122, 9).Value = 1325  # I122: 1400 -> 1325
$ws.Cells.Item(122, 11).Value = 3975  # K122: 4200 -> 3975
$ws.Cells.Item(122, 13).Value = -1525  # M122: -1750 -> -1525
$ws.Cells.Item(132, 8).Value = 1431  # H132: 1167.3636 -> 1431
$ws.Cells.Item(132, 9).Value = 1095.8  # I132: 859 -> 1095.8
$ws.Cells.Item(132, 11).Value = 3287.4  # K132: 2577 -> 3287.4
$ws.Cells.Item(132, 13).Value = -757.3999999999996  # M132: -47 -> -757.3999999999996
